# standard-observations.xlsx update
# - dct:modified (B24) bumped to the date of this regeneration
# - rows 37-86 and 92-104 (the iop:VariableSet concepts) gain the
#   skos:Concept rdf:type alongside iop:VariableSet
# - a handful of those rows' own "modified" date (column U) is bumped too

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Force the cell to keep a literal text value (Excel would otherwise
    # auto-coerce date-shaped strings like "2025-06-13" into a date
    # serial number). Apply a text format while assigning, then clear
    # the formatting again so no extra cell style is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# dct:modified date at the top of the sheet
Set-TextValue $ws.Range("B24") "2025-06-13"

# rdf:type column (D) for the VariableSet concept rows
$variableSetRows = @(37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,92,93,94,95,96,97,98,99,100,101,102,103,104)
foreach ($row in $variableSetRows) {
    $ws.Range("D$row").Value = "skos:Concept, iop:VariableSet"
}

# per-row modified date (column U) bumped for these specific rows
$modifiedDateRows = @(37,38,72,92,93)
foreach ($row in $modifiedDateRows) {
    Set-TextValue $ws.Range("U$row") "2025-06-13"
}
